$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.947.55'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '1.639.23'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  +0.69%  '
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('E6').Value = '  +0.49%  '
$ws.Range('E7').Value = '  +0.67%  '
$ws.Range('E8').Value = '  -0.57%  '
$ws.Range('E9').Value = '  +0.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.64'
$ws.Range('E10').Value = '  -0.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0794'
$ws.Range('E11').Value = '  +0.70%  '
$ws.Range('D12').Value = '1.866.05'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('E13').Value = '  -0.22%  '
$ws.Range('D14').Value = '1.648.14'
$ws.Range('E14').Value = '  +0.87%  '
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.67'
$ws.Range('E17').Value = '  -0.93%  '
$ws.Range('D18').Value = '25.961.27'
$ws.Range('E18').Value = '  -0.22%  '
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '194.18'
$ws.Range('E20').Value = '  +0.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.37'
$ws.Range('E21').Value = '  -1.59%  '
$ws.Range('E22').Value = '  -0.85%  '
$ws.Range('E23').Value = '  -1.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '144.17'
$ws.Range('E24').Value = '  +1.13%  '
$ws.Range('B25').Value = 'BinanceUSD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.69%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.78'
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('E27').Value = '  +1.86%  '
$ws.Range('E28').Value = '  -0.88%  '
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('E31').Value = '  +0.31%  '
$ws.Range('E32').Value = '  -1.27%  '
$ws.Range('E33').Value = '  -0.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.54'
$ws.Range('E34').Value = '  -3.03%  '
$ws.Range('E35').Value = '  +1.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.904'
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('D37').Value = '1.139.54'
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('E38').Value = '  -0.98%  '
$ws.Range('E39').Value = '  -1.88%  '
$ws.Range('E40').Value = '  +0.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '99.30'
$ws.Range('E41').Value = '  -1.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.800'
$ws.Range('E42').Value = '  +0.97%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.43'
$ws.Range('E43').Value = '  -2.98%  '
$ws.Range('D44').Value = '1.776.17'
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('E45').Value = '  +9.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.57'
$ws.Range('E46').Value = '  +1.23%  '
$ws.Range('E47').Value = '  +2.67%  '
$ws.Range('E48').Value = '  -1.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.66'
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('E50').Value = '  -0.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0963'
$ws.Range('E51').Value = '  -1.03%  '
